$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: refreshed "updated at" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 00:05"

# --- Estados Unidos (row 4): refreshed totals ---
$ws.Cells.Item(4, 2).Value = 1526042
$ws.Cells.Item(4, 3).Value = 18269
$ws.Cells.Item(4, 4).Value = 343204
$ws.Cells.Item(4, 5).Value = 1091910
$ws.Cells.Item(4, 7).Value = 815
$ws.Cells.Item(4, 8).Value = 90928

# --- Brasil (row 8): refreshed totals ---
$ws.Cells.Item(8, 2).Value = 239483
$ws.Cells.Item(8, 3).Value = 6341
$ws.Cells.Item(8, 5).Value = 133749
$ws.Cells.Item(8, 7).Value = 429
$ws.Cells.Item(8, 8).Value = 16062

# --- Principado de Andorra (row 116): refreshed totals ---
$ws.Cells.Item(116, 4).Value = 617
$ws.Cells.Item(116, 5).Value = 93

# --- Rows 186/187: Antigua y Barbuda / Botsuana swap rank (tie broken
#     differently after refresh) plus refreshed totals ---
$ws.Cells.Item(186, 1).Value = "Botsuana"
$ws.Cells.Item(186, 3).Value = 1
$ws.Cells.Item(186, 4).Value = 17
$ws.Cells.Item(186, 5).Value = 7
$ws.Cells.Item(186, 8).Value = 1

$ws.Cells.Item(187, 1).Value = "Antigua y Barbuda"
$ws.Cells.Item(187, 2).Value = 25
$ws.Cells.Item(187, 4).Value = 19
$ws.Cells.Item(187, 5).Value = 3
$ws.Cells.Item(187, 8).Value = 3

# --- Rows 195/197: Nueva Caledonia / Santa Lucia swap rank (tied totals,
#     no value changes) ---
$ws.Cells.Item(195, 1).Value = "Santa Lucia"
$ws.Cells.Item(197, 1).Value = "Nueva Caledonia"

# --- Rows 215/216: San Bartolome / Sahara Occidental swap rank (tied
#     totals, no value changes) ---
$ws.Cells.Item(215, 1).Value = "Sahara Occidental"
$ws.Cells.Item(216, 1).Value = "San Bartolome"
